# 교수등록.xlsx :: 단체등록엑셀 - rework the Sheet1 header row
#   - drop the 교수번호/비밀번호 columns, add 학과코드
#   - reorder & relabel the remaining header cells
#   - shrink the sheet from 13 columns (A:M) down to 11 columns (A:K)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the two trailing columns (L, M) entirely so the used range / dimension
# shrinks from A1:M1 down to A1:K1, just like row shifting in the UI.
$ws.Range("L1:M1").EntireColumn.Delete() | Out-Null

# Re-type the header row, in its new left-to-right order.
$ws.Range("A1").Value = "주민번호"
$ws.Range("B1").Value = "학과코드"
$ws.Range("C1").Value = "이름"
$ws.Range("D1").Value = "영문이름"
$ws.Range("E1").Value = "이메일"
$ws.Range("F1").Value = "우편번호"
$ws.Range("G1").Value = "주소"
$ws.Range("H1").Value = "상세 주소"
$ws.Range("I1").Value = "핸드폰 번호"
$ws.Range("J1").Value = "집 전화번호"
$ws.Range("K1").Value = "교수 전화번호"

# Re-fit the column widths to the (now shorter) header captions.
$ws.Columns.Item(1).ColumnWidth = 8.285714285714286
$ws.Columns.Item(3).ColumnWidth = 4.571428571428571
$ws.Columns.Item(5).ColumnWidth = 6.428571428571429
$ws.Columns.Item(7).ColumnWidth = 4.571428571428571
$ws.Columns.Item(8).ColumnWidth = 8.857142857142858
$ws.Columns.Item(9).ColumnWidth = 10.857142857142858
$ws.Columns.Item(10).ColumnWidth = 10.857142857142858
$ws.Columns.Item(11).ColumnWidth = 13.0

# Land the selection where the cursor naturally ends up after typing the
# last header value (one cell past the new last column, K).
$ws.Range("L1").Select() | Out-Null
